$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 234, pushing the existing rows 234-273 down
# to 235-274 (this also carries the date-column number format down with
# the shifted cells).
$ws.Rows.Item(234).Insert()

# Populate the newly-inserted row 234 with a new data record (same
# market/category/etc. as the surrounding rows, new date + volume/price
# figures, origin "Región del Maule").
$ws.Cells.Item(234, 1).Value = 5
$ws.Cells.Item(234, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(234, 3).Value = "Maule"
$ws.Cells.Item(234, 4).Value = 44694
$ws.Cells.Item(234, 5).Value = 7
$ws.Cells.Item(234, 6).Value = 100112006
$ws.Cells.Item(234, 7).Value = "Repollo"
$ws.Cells.Item(234, 8).Value = "Crespo record"
$ws.Cells.Item(234, 9).Value = "Primera"
$ws.Cells.Item(234, 10).Value = 3000
$ws.Cells.Item(234, 11).Value = 1000
$ws.Cells.Item(234, 12).Value = 1000
$ws.Cells.Item(234, 13).Value = 1000
$ws.Cells.Item(234, 14).Value = "`$/unidad"
$ws.Cells.Item(234, 15).Value = "Región del Maule"
$ws.Cells.Item(234, 16).Value = 1000
$ws.Cells.Item(234, 17).Value = 1
$ws.Cells.Item(234, 18).Value = "Hortaliza"
